$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "17:20 03-Dec-23"
$ws.Range("C10").Value = "Ẩn danh"
$ws.Range("D10").Value = "fsad"
